# Update the date line and every "a OP b = c" answer in the table to the
# new values from the target revision. Each old string is unique within the
# document, so a simple whole-document Find/Replace (whole word match) for
# each pair is safe and unambiguous.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-30 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-07-31 Monday", 2) | Out-Null
$d.Content.Find.Execute("61+19=80", $true, $true, $false, $false, $false, $true, 1, $false, "91-31=60", 2) | Out-Null
$d.Content.Find.Execute("12+37=49", $true, $true, $false, $false, $false, $true, 1, $false, "94-81=13", 2) | Out-Null
$d.Content.Find.Execute("20+73=93", $true, $true, $false, $false, $false, $true, 1, $false, "18+9=27", 2) | Out-Null
$d.Content.Find.Execute("31-25=6", $true, $true, $false, $false, $false, $true, 1, $false, "91-59=32", 2) | Out-Null
$d.Content.Find.Execute("59-18=41", $true, $true, $false, $false, $false, $true, 1, $false, "32-27=5", 2) | Out-Null
$d.Content.Find.Execute("10+89=99", $true, $true, $false, $false, $false, $true, 1, $false, "6+42=48", 2) | Out-Null
$d.Content.Find.Execute("3+31=34", $true, $true, $false, $false, $false, $true, 1, $false, "38-28=10", 2) | Out-Null
$d.Content.Find.Execute("72-13=59", $true, $true, $false, $false, $false, $true, 1, $false, "9+15=24", 2) | Out-Null
$d.Content.Find.Execute("44+25=69", $true, $true, $false, $false, $false, $true, 1, $false, "64-62=2", 2) | Out-Null
$d.Content.Find.Execute("88-73=15", $true, $true, $false, $false, $false, $true, 1, $false, "91+6=97", 2) | Out-Null
$d.Content.Find.Execute("92-65=27", $true, $true, $false, $false, $false, $true, 1, $false, "89-83=6", 2) | Out-Null
$d.Content.Find.Execute("80+8=88", $true, $true, $false, $false, $false, $true, 1, $false, "69+26=95", 2) | Out-Null
$d.Content.Find.Execute("94-83=11", $true, $true, $false, $false, $false, $true, 1, $false, "78-10=68", 2) | Out-Null
$d.Content.Find.Execute("91-43=48", $true, $true, $false, $false, $false, $true, 1, $false, "90-25=65", 2) | Out-Null
$d.Content.Find.Execute("69-68=1", $true, $true, $false, $false, $false, $true, 1, $false, "4+37=41", 2) | Out-Null
$d.Content.Find.Execute("34-30=4", $true, $true, $false, $false, $false, $true, 1, $false, "89-86=3", 2) | Out-Null
$d.Content.Find.Execute("41-36=5", $true, $true, $false, $false, $false, $true, 1, $false, "10+49=59", 2) | Out-Null
$d.Content.Find.Execute("55+7=62", $true, $true, $false, $false, $false, $true, 1, $false, "50+12=62", 2) | Out-Null
$d.Content.Find.Execute("94+0=94", $true, $true, $false, $false, $false, $true, 1, $false, "69+10=79", 2) | Out-Null
$d.Content.Find.Execute("81+0=81", $true, $true, $false, $false, $false, $true, 1, $false, "12+39=51", 2) | Out-Null
$d.Content.Find.Execute("25+26=51", $true, $true, $false, $false, $false, $true, 1, $false, "73+13=86", 2) | Out-Null
$d.Content.Find.Execute("99-42=57", $true, $true, $false, $false, $false, $true, 1, $false, "54-3=51", 2) | Out-Null
$d.Content.Find.Execute("52+25=77", $true, $true, $false, $false, $false, $true, 1, $false, "33-4=29", 2) | Out-Null
$d.Content.Find.Execute("10+77=87", $true, $true, $false, $false, $false, $true, 1, $false, "73-64=9", 2) | Out-Null
$d.Content.Find.Execute("41+51=92", $true, $true, $false, $false, $false, $true, 1, $false, "44+2=46", 2) | Out-Null
$d.Content.Find.Execute("87-23=64", $true, $true, $false, $false, $false, $true, 1, $false, "40+35=75", 2) | Out-Null
$d.Content.Find.Execute("20-2=18", $true, $true, $false, $false, $false, $true, 1, $false, "19+22=41", 2) | Out-Null
$d.Content.Find.Execute("0+42=42", $true, $true, $false, $false, $false, $true, 1, $false, "49-47=2", 2) | Out-Null
$d.Content.Find.Execute("14+33=47", $true, $true, $false, $false, $false, $true, 1, $false, "90-60=30", 2) | Out-Null
$d.Content.Find.Execute("2+26=28", $true, $true, $false, $false, $false, $true, 1, $false, "8+73=81", 2) | Out-Null
$d.Content.Find.Execute("65-46=19", $true, $true, $false, $false, $false, $true, 1, $false, "44-18=26", 2) | Out-Null
$d.Content.Find.Execute("3+24=27", $true, $true, $false, $false, $false, $true, 1, $false, "18+24=42", 2) | Out-Null
$d.Content.Find.Execute("98-23=75", $true, $true, $false, $false, $false, $true, 1, $false, "84+8=92", 2) | Out-Null
$d.Content.Find.Execute("17+78=95", $true, $true, $false, $false, $false, $true, 1, $false, "88-22=66", 2) | Out-Null
$d.Content.Find.Execute("85-3=82", $true, $true, $false, $false, $false, $true, 1, $false, "6+52=58", 2) | Out-Null
$d.Content.Find.Execute("5+29=34", $true, $true, $false, $false, $false, $true, 1, $false, "38-6=32", 2) | Out-Null
$d.Content.Find.Execute("9+65=74", $true, $true, $false, $false, $false, $true, 1, $false, "39+44=83", 2) | Out-Null
$d.Content.Find.Execute("33-9=24", $true, $true, $false, $false, $false, $true, 1, $false, "78-10=68", 2) | Out-Null
$d.Content.Find.Execute("48+40=88", $true, $true, $false, $false, $false, $true, 1, $false, "72-10=62", 2) | Out-Null
$d.Content.Find.Execute("5+28=33", $true, $true, $false, $false, $false, $true, 1, $false, "78-14=64", 2) | Out-Null
$d.Content.Find.Execute("69-43=26", $true, $true, $false, $false, $false, $true, 1, $false, "40+26=66", 2) | Out-Null
$d.Content.Find.Execute("49+29=78", $true, $true, $false, $false, $false, $true, 1, $false, "45+48=93", 2) | Out-Null
$d.Content.Find.Execute("71-38=33", $true, $true, $false, $false, $false, $true, 1, $false, "65-14=51", 2) | Out-Null
$d.Content.Find.Execute("61-20=41", $true, $true, $false, $false, $false, $true, 1, $false, "55+21=76", 2) | Out-Null
$d.Content.Find.Execute("63-44=19", $true, $true, $false, $false, $false, $true, 1, $false, "50-40=10", 2) | Out-Null
$d.Content.Find.Execute("43-34=9", $true, $true, $false, $false, $false, $true, 1, $false, "92-39=53", 2) | Out-Null
$d.Content.Find.Execute("10+76=86", $true, $true, $false, $false, $false, $true, 1, $false, "93+3=96", 2) | Out-Null
$d.Content.Find.Execute("6+37=43", $true, $true, $false, $false, $false, $true, 1, $false, "80+2=82", 2) | Out-Null
$d.Content.Find.Execute("3+81=84", $true, $true, $false, $false, $false, $true, 1, $false, "53+12=65", 2) | Out-Null
$d.Content.Find.Execute("7+0=7", $true, $true, $false, $false, $false, $true, 1, $false, "20-19=1", 2) | Out-Null
$d.Content.Find.Execute("30+62=92", $true, $true, $false, $false, $false, $true, 1, $false, "27+46=73", 2) | Out-Null
$d.Content.Find.Execute("21-7=14", $true, $true, $false, $false, $false, $true, 1, $false, "84-49=35", 2) | Out-Null
$d.Content.Find.Execute("87-39=48", $true, $true, $false, $false, $false, $true, 1, $false, "39-39=0", 2) | Out-Null
$d.Content.Find.Execute("61-34=27", $true, $true, $false, $false, $false, $true, 1, $false, "2+45=47", 2) | Out-Null
$d.Content.Find.Execute("35+54=89", $true, $true, $false, $false, $false, $true, 1, $false, "93-38=55", 2) | Out-Null
$d.Content.Find.Execute("67-22=45", $true, $true, $false, $false, $false, $true, 1, $false, "19+6=25", 2) | Out-Null
$d.Content.Find.Execute("77-42=35", $true, $true, $false, $false, $false, $true, 1, $false, "89-54=35", 2) | Out-Null
$d.Content.Find.Execute("34+11=45", $true, $true, $false, $false, $false, $true, 1, $false, "59-50=9", 2) | Out-Null
$d.Content.Find.Execute("32+42=74", $true, $true, $false, $false, $false, $true, 1, $false, "16+31=47", 2) | Out-Null
$d.Content.Find.Execute("60-9=51", $true, $true, $false, $false, $false, $true, 1, $false, "21+74=95", 2) | Out-Null
$d.Content.Find.Execute("80+4=84", $true, $true, $false, $false, $false, $true, 1, $false, "9+51=60", 2) | Out-Null
$d.Content.Find.Execute("70-56=14", $true, $true, $false, $false, $false, $true, 1, $false, "68-56=12", 2) | Out-Null
$d.Content.Find.Execute("70-44=26", $true, $true, $false, $false, $false, $true, 1, $false, "46-39=7", 2) | Out-Null
$d.Content.Find.Execute("36+55=91", $true, $true, $false, $false, $false, $true, 1, $false, "67+27=94", 2) | Out-Null
$d.Content.Find.Execute("47+51=98", $true, $true, $false, $false, $false, $true, 1, $false, "19+76=95", 2) | Out-Null
$d.Content.Find.Execute("24+58=82", $true, $true, $false, $false, $false, $true, 1, $false, "51+22=73", 2) | Out-Null
$d.Content.Find.Execute("23+29=52", $true, $true, $false, $false, $false, $true, 1, $false, "49+35=84", 2) | Out-Null
$d.Content.Find.Execute("72-11=61", $true, $true, $false, $false, $false, $true, 1, $false, "2+3=5", 2) | Out-Null
$d.Content.Find.Execute("14-0=14", $true, $true, $false, $false, $false, $true, 1, $false, "66-63=3", 2) | Out-Null
$d.Content.Find.Execute("42-35=7", $true, $true, $false, $false, $false, $true, 1, $false, "97-74=23", 2) | Out-Null
$d.Content.Find.Execute("18+54=72", $true, $true, $false, $false, $false, $true, 1, $false, "23+0=23", 2) | Out-Null
$d.Content.Find.Execute("78+11=89", $true, $true, $false, $false, $false, $true, 1, $false, "30-13=17", 2) | Out-Null
$d.Content.Find.Execute("86-53=33", $true, $true, $false, $false, $false, $true, 1, $false, "37-2=35", 2) | Out-Null
$d.Content.Find.Execute("1+64=65", $true, $true, $false, $false, $false, $true, 1, $false, "76-67=9", 2) | Out-Null
$d.Content.Find.Execute("46+30=76", $true, $true, $false, $false, $false, $true, 1, $false, "41-16=25", 2) | Out-Null
$d.Content.Find.Execute("96-90=6", $true, $true, $false, $false, $false, $true, 1, $false, "40-9=31", 2) | Out-Null
$d.Content.Find.Execute("10+71=81", $true, $true, $false, $false, $false, $true, 1, $false, "75-10=65", 2) | Out-Null
$d.Content.Find.Execute("33+25=58", $true, $true, $false, $false, $false, $true, 1, $false, "36-30=6", 2) | Out-Null
$d.Content.Find.Execute("59-27=32", $true, $true, $false, $false, $false, $true, 1, $false, "41+42=83", 2) | Out-Null
$d.Content.Find.Execute("97-54=43", $true, $true, $false, $false, $false, $true, 1, $false, "66+28=94", 2) | Out-Null
$d.Content.Find.Execute("14+49=63", $true, $true, $false, $false, $false, $true, 1, $false, "0+87=87", 2) | Out-Null
$d.Content.Find.Execute("77+7=84", $true, $true, $false, $false, $false, $true, 1, $false, "93-70=23", 2) | Out-Null
$d.Content.Find.Execute("95-74=21", $true, $true, $false, $false, $false, $true, 1, $false, "87+2=89", 2) | Out-Null
$d.Content.Find.Execute("34-6=28", $true, $true, $false, $false, $false, $true, 1, $false, "49+38=87", 2) | Out-Null
$d.Content.Find.Execute("48+0=48", $true, $true, $false, $false, $false, $true, 1, $false, "49-43=6", 2) | Out-Null
$d.Content.Find.Execute("83-35=48", $true, $true, $false, $false, $false, $true, 1, $false, "37+46=83", 2) | Out-Null
$d.Content.Find.Execute("30-4=26", $true, $true, $false, $false, $false, $true, 1, $false, "33+63=96", 2) | Out-Null
$d.Content.Find.Execute("16+46=62", $true, $true, $false, $false, $false, $true, 1, $false, "7+6=13", 2) | Out-Null
$d.Content.Find.Execute("84-68=16", $true, $true, $false, $false, $false, $true, 1, $false, "87-4=83", 2) | Out-Null
$d.Content.Find.Execute("20-5=15", $true, $true, $false, $false, $false, $true, 1, $false, "85-13=72", 2) | Out-Null
$d.Content.Find.Execute("37-19=18", $true, $true, $false, $false, $false, $true, 1, $false, "52-44=8", 2) | Out-Null
$d.Content.Find.Execute("43-2=41", $true, $true, $false, $false, $false, $true, 1, $false, "12-5=7", 2) | Out-Null
$d.Content.Find.Execute("68-22=46", $true, $true, $false, $false, $false, $true, 1, $false, "42-32=10", 2) | Out-Null
$d.Content.Find.Execute("30+37=67", $true, $true, $false, $false, $false, $true, 1, $false, "28+21=49", 2) | Out-Null
$d.Content.Find.Execute("95-47=48", $true, $true, $false, $false, $false, $true, 1, $false, "25+33=58", 2) | Out-Null
$d.Content.Find.Execute("1+10=11", $true, $true, $false, $false, $false, $true, 1, $false, "78-20=58", 2) | Out-Null
$d.Content.Find.Execute("60+30=90", $true, $true, $false, $false, $false, $true, 1, $false, "22-9=13", 2) | Out-Null
$d.Content.Find.Execute("44+52=96", $true, $true, $false, $false, $false, $true, 1, $false, "54+27=81", 2) | Out-Null
$d.Content.Find.Execute("23+62=85", $true, $true, $false, $false, $false, $true, 1, $false, "81-7=74", 2) | Out-Null
$d.Content.Find.Execute("4+49=53", $true, $true, $false, $false, $false, $true, 1, $false, "82-40=42", 2) | Out-Null
